# Tonberry_Profits: refresh market-board price snapshots + recomputed leve profit columns.
# Values below were taken verbatim from the updated data pull (see commit message).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2500
$ws.Range("J32").Value = 2500
$ws.Range("L32").Value = 2500
$ws.Range("N32").Value = -3152
$ws.Range("H87").Value = 25000
$ws.Range("J87").Value = 25000
$ws.Range("L87").Value = 25000
$ws.Range("N87").Value = -27496
$ws.Range("H90").Value = 25000
$ws.Range("J90").Value = 25000
$ws.Range("L90").Value = 75000
$ws.Range("N90").Value = -87480
$ws.Range("H129").Value = 873.05457
$ws.Range("I129").Value = 699.4
$ws.Range("K129").Value = 2098.2
$ws.Range("M129").Value = 2901.8
$ws.Range("H135").Value = 349.82608
$ws.Range("I135").Value = 373.7619
$ws.Range("J135").Value = 98.5
$ws.Range("K135").Value = 3363.8571
$ws.Range("L135").Value = 886.5
$ws.Range("M135").Value = -828.8571000000002
$ws.Range("N135").Value = -5956.5
$ws.Range("H137").Value = 1950.3572
$ws.Range("I137").Value = 1625.5
$ws.Range("J137").Value = 2194
$ws.Range("K137").Value = 4876.5
$ws.Range("L137").Value = 6582
$ws.Range("M137").Value = -2326.5
$ws.Range("N137").Value = -11682

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1651.3334
$ws.Range("I97").Value = 1227.25
$ws.Range("J97").Value = 2499.5
$ws.Range("K97").Value = 1227.25
$ws.Range("L97").Value = 2499.5
$ws.Range("M97").Value = -731.25
$ws.Range("N97").Value = -3491.5
$ws.Range("H110").Value = 1555.5
$ws.Range("I110").Value = 1187.7826
$ws.Range("K110").Value = 1187.7826
$ws.Range("M110").Value = 857.2174
$ws.Range("H132").Value = 1428
$ws.Range("I132").Value = 975.5897
$ws.Range("K132").Value = 2926.7691
$ws.Range("M132").Value = -396.7691

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 17348
$ws.Range("J80").Value = 20717.6
$ws.Range("L80").Value = 20717.6
$ws.Range("N80").Value = -22713.6
$ws.Range("H83").Value = 17348
$ws.Range("J83").Value = 20717.6
$ws.Range("L83").Value = 103588
$ws.Range("N83").Value = -113572
$ws.Range("H86").Value = 128286.75
$ws.Range("I86").Value = 6156.4
$ws.Range("J86").Value = 183800.55
$ws.Range("K86").Value = 6156.4
$ws.Range("L86").Value = 183800.55
$ws.Range("M86").Value = -5033.4
$ws.Range("N86").Value = -186046.55
$ws.Range("H89").Value = 128286.75
$ws.Range("I89").Value = 6156.4
$ws.Range("J89").Value = 183800.55
$ws.Range("K89").Value = 30782
$ws.Range("L89").Value = 919002.75
$ws.Range("M89").Value = -25166
$ws.Range("N89").Value = -930234.75
$ws.Range("H94").Value = 1034.5
$ws.Range("I94").Value = 1034.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1034.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -583.5
$ws.Range("N94").ClearContents()
$ws.Range("H105").Value = 2541
$ws.Range("I105").Value = 2541
$ws.Range("K105").Value = 2541
$ws.Range("M105").Value = -794
$ws.Range("H107").Value = 1074.9286
$ws.Range("I107").Value = 1080.6154
$ws.Range("K107").Value = 1080.6154
$ws.Range("M107").Value = 839.3846000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2520.625
$ws.Range("J31").Value = 2550.1765
$ws.Range("L31").Value = 2550.1765
$ws.Range("N31").Value = -3140.1765
$ws.Range("H34").Value = 2520.625
$ws.Range("J34").Value = 2550.1765
$ws.Range("L34").Value = 2550.1765
$ws.Range("N34").Value = -2954.1765
$ws.Range("H107").Value = 577.3
$ws.Range("J107").Value = 760
$ws.Range("L107").Value = 760
$ws.Range("N107").Value = -4600

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 379.53845
$ws.Range("I7").Value = 190.28572
$ws.Range("J7").Value = 600.3333
$ws.Range("K7").Value = 570.85716
$ws.Range("L7").Value = 1800.9999
$ws.Range("M7").Value = -458.85716
$ws.Range("N7").Value = -2024.9999
$ws.Range("H131").Value = 12398.246
$ws.Range("J131").Value = 13513.048
$ws.Range("L131").Value = 40539.144
$ws.Range("N131").Value = -50619.144
$ws.Range("H136").Value = 1356.3334
$ws.Range("I136").Value = 1356.3334
$ws.Range("K136").Value = 4069.0002
$ws.Range("M136").Value = 1030.9998
$ws.Range("H138").Value = 3663.125
$ws.Range("I138").Value = 2925.8333
$ws.Range("K138").Value = 8777.499899999999
$ws.Range("M138").Value = -3637.499899999999
$ws.Range("H141").Value = 5188.4
$ws.Range("I141").Value = 5188.4
$ws.Range("K141").Value = 15565.2
$ws.Range("M141").Value = -10385.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 932.73334
$ws.Range("I113").Value = 705.3333
$ws.Range("K113").Value = 705.3333
$ws.Range("M113").Value = 1464.6667
$ws.Range("H132").Value = 1071147.9
$ws.Range("I132").Value = 1540857.5
$ws.Range("J132").Value = 3625.9092
$ws.Range("K132").Value = 4622572.5
$ws.Range("L132").Value = 10877.7276
$ws.Range("M132").Value = -4620042.5
$ws.Range("N132").Value = -15937.7276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3207
$ws.Range("J22").Value = 2333.3333
$ws.Range("L22").Value = 2333.3333
$ws.Range("N22").Value = -2923.3333
$ws.Range("H27").Value = 3207
$ws.Range("J27").Value = 2333.3333
$ws.Range("L27").Value = 2333.3333
$ws.Range("N27").Value = -2547.3333
$ws.Range("H40").Value = 2855.5625
$ws.Range("I40").Value = 2254.0908
$ws.Range("K40").Value = 2254.0908
$ws.Range("M40").Value = -2118.0908
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H47").Value = 25000
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 25000
$ws.Range("K47").Value = 0
$ws.Range("L47").Value = 25000
$ws.Range("M47").ClearContents()
$ws.Range("N47").Value = -25980
$ws.Range("H52").Value = 25000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 25000
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 25000
$ws.Range("M52").ClearContents()
$ws.Range("N52").Value = -25466
$ws.Range("H55").Value = 297.38235
$ws.Range("I55").Value = 262.24
$ws.Range("J55").Value = 395
$ws.Range("K55").Value = 262.24
$ws.Range("L55").Value = 395
$ws.Range("M55").Value = -89.24000000000001
$ws.Range("N55").Value = -741
$ws.Range("H61").Value = 4249.75
$ws.Range("I61").Value = 3499.5
$ws.Range("K61").Value = 3499.5
$ws.Range("M61").Value = -3297.5
$ws.Range("H113").Value = 4249.75
$ws.Range("I113").Value = 3499.5
$ws.Range("K113").Value = 3499.5
$ws.Range("M113").Value = -1329.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 618.3333
$ws.Range("I113").Value = 345.3
$ws.Range("K113").Value = 1035.9
$ws.Range("M113").Value = 1134.1
$ws.Range("H132").Value = 1700.25
$ws.Range("I132").Value = 1087.6207
$ws.Range("K132").Value = 3262.8621
$ws.Range("M132").Value = -732.8620999999998
$ws.Range("H136").Value = 23150574
$ws.Range("J136").Value = 3500
$ws.Range("L136").Value = 10500
$ws.Range("N136").Value = -15600

Write-Host "Applied 181 value updates and 4 clears across 8 sheets."
